$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 data (ENZt2r enzyme entry)
$ws.Range("A8").Value = "ENZt2r"
$ws.Range("C8").Value = "enz1[c] <==>"
$ws.Range("E8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1

# Column widths for T (20) and U (21), matching the new data region
# (values chosen so the engine's pixel-grid rounding lands closest to the
# author's original 12.28515625 / 12.85546875 character widths)
$ws.Columns.Item(20).ColumnWidth = 11.5
$ws.Columns.Item(21).ColumnWidth = 12

# Update the selection to cover the new used range
$null = $ws.Range("A1:U8").Select()
